# ------------------------------------------------------------------
# Horarios Linea 141 update (11/01/2026 run @ 11:52:01)
# Regenerates the schedule sheet ("LP1912") cell-by-cell to match the
# latest scrape, refreshes the "Ultima actualizacion" / "Total filas"
# headers on all three sheets, and appends the newly scraped rows.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$sheetLP1912    = $wb.Worksheets.Item("LP1912")
$sheetLP1912215 = $wb.Worksheets.Item("LP1912-215")
$sheet6203      = $wb.Worksheets.Item("6203-6173")

$newTimestamp = "11:52:01"

# --- Header block on the main data sheet -------------------------
$sheetLP1912.Range("A2").Value = "Última actualización: $newTimestamp"
$sheetLP1912.Range("A3").Value = "Total filas: 149"

# --- Header timestamp on the other two sheets (row counts unchanged)
$sheetLP1912215.Range("A2").Value = "Última actualización: $newTimestamp"
$sheet6203.Range("A2").Value = "Última actualización: $newTimestamp"

# --- Cell-level corrections on existing rows 23-150 ---------------
# Each tuple is (row, column letter, new value)
$changes = @(
    ,@(23, "A", "05:57:13")
    ,@(23, "C", "23_HERNANDEZ")
    ,@(23, "D", 84)
    ,@(24, "A", "06:17:28")
    ,@(24, "C", "16_SANTA ANA")
    ,@(24, "D", 64)
    ,@(33, "A", "06:46:50")
    ,@(33, "C", "16_SANTA ANA")
    ,@(33, "D", 74)
    ,@(34, "A", "06:17:28")
    ,@(34, "C", "17_ROMERO")
    ,@(34, "D", 103)
    ,@(99, "A", "11:11:33")
    ,@(99, "C", "215C_EL PATO")
    ,@(99, "D", 20)
    ,@(100, "A", "10:05:51")
    ,@(100, "C", "16_SANTA ANA")
    ,@(100, "D", 86)
    ,@(106, "A", "11:47:17")
    ,@(106, "C", "23_HERNANDEZ")
    ,@(106, "D", 5)
    ,@(107, "A", "11:52:01")
    ,@(107, "C", "15X38_ABASTO")
    ,@(107, "D", 0)
    ,@(108, "A", "10:05:51")
    ,@(108, "B", "11:52")
    ,@(108, "D", 107)
    ,@(109, "A", "10:50:41")
    ,@(109, "C", "225_GOMEZ")
    ,@(109, "D", 63)
    ,@(110, "A", "10:37:52")
    ,@(110, "B", "11:53")
    ,@(110, "D", 76)
    ,@(111, "A", "11:52:01")
    ,@(111, "B", "11:54")
    ,@(111, "C", "225_GOMEZ")
    ,@(111, "D", 2)
    ,@(112, "A", "10:50:41")
    ,@(112, "B", "11:54")
    ,@(112, "C", "23_HERNANDEZ")
    ,@(112, "D", 64)
    ,@(113, "A", "11:34:59")
    ,@(113, "B", "11:57")
    ,@(113, "C", "17_ROMERO")
    ,@(113, "D", 23)
    ,@(114, "A", "10:05:51")
    ,@(114, "B", "11:58")
    ,@(114, "C", "17_ROMERO")
    ,@(114, "D", 113)
    ,@(115, "A", "10:37:52")
    ,@(115, "B", "12:05")
    ,@(115, "C", "11_ETCHEVERRY")
    ,@(115, "D", 88)
    ,@(116, "A", "11:47:17")
    ,@(116, "B", "12:06")
    ,@(116, "C", "11_ETCHEVERRY")
    ,@(116, "D", 19)
    ,@(117, "A", "11:34:59")
    ,@(117, "B", "12:09")
    ,@(117, "C", "16_P MOR-SANTA ANA")
    ,@(117, "D", 35)
    ,@(118, "A", "11:34:59")
    ,@(118, "B", "12:09")
    ,@(118, "C", "15_ABASTO")
    ,@(118, "D", 35)
    ,@(119, "B", "12:10")
    ,@(119, "C", "15_ABASTO")
    ,@(119, "D", 93)
    ,@(120, "A", "10:37:52")
    ,@(120, "B", "12:10")
    ,@(120, "C", "16_P MOR-SANTA ANA")
    ,@(120, "D", 93)
    ,@(121, "B", "12:16")
    ,@(121, "C", "10_OLMOS")
    ,@(121, "D", 99)
    ,@(122, "B", "12:17")
    ,@(122, "C", "10_OLMOS")
    ,@(122, "D", 66)
    ,@(123, "A", "10:37:52")
    ,@(123, "B", "12:21")
    ,@(123, "C", "215C_EL PATO")
    ,@(123, "D", 104)
    ,@(124, "A", "11:11:33")
    ,@(124, "B", "12:22")
    ,@(124, "C", "215C_EL PATO")
    ,@(124, "D", 71)
    ,@(125, "B", "12:32")
    ,@(125, "C", "23_HERNANDEZ")
    ,@(125, "D", 45)
    ,@(126, "A", "10:37:52")
    ,@(126, "B", "12:32")
    ,@(126, "C", "14_ABASTO")
    ,@(126, "D", 115)
    ,@(127, "A", "11:47:17")
    ,@(127, "B", "12:33")
    ,@(127, "C", "14_ABASTO")
    ,@(127, "D", 46)
    ,@(128, "A", "11:34:59")
    ,@(128, "B", "12:33")
    ,@(128, "C", "15_ABASTO")
    ,@(128, "D", 59)
    ,@(129, "A", "10:37:52")
    ,@(129, "B", "12:34")
    ,@(129, "C", "15_ABASTO")
    ,@(129, "D", 117)
    ,@(130, "A", "11:11:33")
    ,@(130, "B", "12:35")
    ,@(130, "C", "23_HERNANDEZ")
    ,@(130, "D", 84)
    ,@(131, "B", "12:35")
    ,@(131, "C", "27_EL RETIRO")
    ,@(131, "D", 61)
    ,@(132, "A", "10:50:41")
    ,@(132, "B", "12:36")
    ,@(132, "D", 106)
    ,@(133, "B", "12:36")
    ,@(133, "C", "23_HERNANDEZ")
    ,@(133, "D", 62)
    ,@(134, "A", "11:47:17")
    ,@(134, "B", "12:37")
    ,@(134, "C", "27_EL RETIRO")
    ,@(134, "D", 50)
    ,@(135, "A", "11:52:01")
    ,@(135, "B", "12:37")
    ,@(135, "C", "23_HERNANDEZ")
    ,@(135, "D", 45)
    ,@(136, "A", "11:34:59")
    ,@(136, "B", "12:47")
    ,@(136, "D", 73)
    ,@(137, "A", "11:34:59")
    ,@(137, "B", "12:47")
    ,@(137, "D", 73)
    ,@(138, "A", "11:34:59")
    ,@(138, "B", "12:47")
    ,@(138, "D", 73)
    ,@(139, "A", "11:47:17")
    ,@(139, "B", "12:48")
    ,@(139, "C", "14_ABASTO")
    ,@(139, "D", 61)
    ,@(140, "A", "11:11:33")
    ,@(140, "B", "12:48")
    ,@(140, "C", "15X38_ABASTO")
    ,@(140, "D", 97)
    ,@(141, "A", "10:50:41")
    ,@(141, "B", "12:48")
    ,@(141, "C", "16_SANTA ANA")
    ,@(141, "D", 118)
    ,@(142, "A", "11:11:33")
    ,@(142, "B", "13:02")
    ,@(142, "C", "11_ETCHEVERRY")
    ,@(142, "D", 111)
    ,@(143, "B", "13:03")
    ,@(143, "C", "215C_EL PATO")
    ,@(143, "D", 89)
    ,@(144, "B", "13:03")
    ,@(144, "C", "11_ETCHEVERRY")
    ,@(144, "D", 76)
    ,@(145, "A", "11:47:17")
    ,@(145, "B", "13:04")
    ,@(145, "C", "215C_EL PATO")
    ,@(145, "D", 77)
    ,@(146, "A", "11:34:59")
    ,@(146, "B", "13:12")
    ,@(146, "C", "16_SANTA ANA")
    ,@(146, "D", 98)
    ,@(147, "A", "11:47:17")
    ,@(147, "B", "13:13")
    ,@(147, "C", "16_SANTA ANA")
    ,@(147, "D", 86)
    ,@(148, "A", "11:34:59")
    ,@(148, "B", "13:16")
    ,@(148, "C", "10_OLMOS")
    ,@(148, "D", 102)
    ,@(149, "A", "11:47:17")
    ,@(149, "B", "13:17")
    ,@(149, "C", "10_OLMOS")
    ,@(149, "D", 90)
    ,@(150, "A", "11:34:59")
    ,@(150, "B", "13:24")
    ,@(150, "C", "16_P MOR-SANTA ANA")
    ,@(150, "D", 110)
)

foreach ($change in $changes) {
    $rowNum = $change[0]
    $colLetter = $change[1]
    $newValue = $change[2]
    $sheetLP1912.Range("$colLetter$rowNum").Value = $newValue
}

# --- Newly scraped rows appended at the bottom (151-154) ---------
$newRows = @(
    ,@(151, "11:47:17", "13:25", "16_P MOR-SANTA ANA", 98, "LP1912")
    ,@(152, "11:34:59", "13:32", "215A_EL PATO", 118, "LP1912")
    ,@(153, "11:47:17", "13:33", "215A_EL PATO", 106, "LP1912")
    ,@(154, "11:52:01", "13:47", "225_GOMEZ", 115, "LP1912")
)

foreach ($newRow in $newRows) {
    $rowNum = $newRow[0]
    $sheetLP1912.Range("A$rowNum").Value = $newRow[1]
    $sheetLP1912.Range("B$rowNum").Value = $newRow[2]
    $sheetLP1912.Range("C$rowNum").Value = $newRow[3]
    $sheetLP1912.Range("D$rowNum").Value = $newRow[4]
    $sheetLP1912.Range("E$rowNum").Value = $newRow[5]
}

